$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3304.0557
$ws.Range("I76").Value = 3170.9285
$ws.Range("K76").Value = 3170.9285
$ws.Range("M76").Value = -2855.9285
$ws.Range("H79").Value = 3304.0557
$ws.Range("I79").Value = 3170.9285
$ws.Range("K79").Value = 3170.9285
$ws.Range("M79").Value = -2078.9285
$ws.Range("H132").Value = 45641296
$ws.Range("I132").Value = 83669544
$ws.Range("J132").Value = 7399.9
$ws.Range("K132").Value = 251008632
$ws.Range("L132").Value = 22199.7
$ws.Range("M132").Value = -251006102
$ws.Range("N132").Value = -27259.7
$ws.Range("H138").Value = 2296.9
$ws.Range("I138").Value = 1271.3422
$ws.Range("J138").Value = 3224.7856
$ws.Range("K138").Value = 3814.0266
$ws.Range("L138").Value = 9674.356800000001
$ws.Range("M138").Value = 1325.9734
$ws.Range("N138").Value = -19954.3568
$ws.Range("H141").Value = 1883.7317
$ws.Range("I141").Value = 1575.0857
$ws.Range("J141").Value = 3684.1667
$ws.Range("K141").Value = 4725.257100000001
$ws.Range("L141").Value = 11052.5001
$ws.Range("M141").Value = 454.7428999999993
$ws.Range("N141").Value = -21412.5001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 534.9286
$ws.Range("I2").Value = 504.1613
$ws.Range("J2").Value = 621.63635
$ws.Range("K2").Value = 504.1613
$ws.Range("L2").Value = 621.63635
$ws.Range("M2").Value = -391.1613
$ws.Range("N2").Value = -847.63635
$ws.Range("H32").Value = 10724.429
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10724.429
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10724.429
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -11298.429
$ws.Range("H61").Value = 2084.9167
$ws.Range("I61").Value = 1294.7858
$ws.Range("J61").Value = 3191.1
$ws.Range("K61").Value = 1294.7858
$ws.Range("L61").Value = 3191.1
$ws.Range("M61").Value = -1082.7858
$ws.Range("N61").Value = -3615.1
$ws.Range("H112").Value = 26500
$ws.Range("J112").Value = 26500
$ws.Range("L112").Value = 26500
$ws.Range("N112").Value = -29454
$ws.Range("H116").Value = 534.9286
$ws.Range("I116").Value = 504.1613
$ws.Range("J116").Value = 621.63635
$ws.Range("K116").Value = 504.1613
$ws.Range("L116").Value = 621.63635
$ws.Range("M116").Value = 1789.8387
$ws.Range("N116").Value = -5209.63635
$ws.Range("H136").Value = 2084.9167
$ws.Range("I136").Value = 1294.7858
$ws.Range("J136").Value = 3191.1
$ws.Range("K136").Value = 3884.3574
$ws.Range("L136").Value = 9573.299999999999
$ws.Range("M136").Value = -1334.3574
$ws.Range("N136").Value = -14673.3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 534.9286
$ws.Range("I3").Value = 504.1613
$ws.Range("J3").Value = 621.63635
$ws.Range("K3").Value = 504.1613
$ws.Range("L3").Value = 621.63635
$ws.Range("M3").Value = -390.1613
$ws.Range("N3").Value = -849.63635
$ws.Range("H5").Value = 1822.875
$ws.Range("I5").Value = 1596.6
$ws.Range("J5").Value = 2200
$ws.Range("K5").Value = 1596.6
$ws.Range("L5").Value = 2200
$ws.Range("M5").Value = -1483.6
$ws.Range("N5").Value = -2426
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5732.476
$ws.Range("I132").Value = 6614.857
$ws.Range("J132").Value = 5291.2856
$ws.Range("K132").Value = 19844.571
$ws.Range("L132").Value = 15873.8568
$ws.Range("M132").Value = -17314.571
$ws.Range("N132").Value = -20933.8568
$ws.Range("H134").Value = 10717.538
$ws.Range("I134").Value = 12501.556
$ws.Range("J134").Value = 6703.5
$ws.Range("K134").Value = 37504.66800000001
$ws.Range("L134").Value = 20110.5
$ws.Range("M134").Value = -34969.66800000001
$ws.Range("N134").Value = -25180.5
$ws.Range("H137").Value = 48540
$ws.Range("J137").Value = 48540
$ws.Range("L137").Value = 48540
$ws.Range("N137").Value = -58740
$ws.Range("H139").Value = 40730
$ws.Range("J139").Value = 40730
$ws.Range("L139").Value = 40730
$ws.Range("N139").Value = -51010
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = ""
$ws.Range("H75").Value = 3477.75
$ws.Range("J75").Value = 3477.75
$ws.Range("L75").Value = 10433.25
$ws.Range("N75").Value = -12429.25
$ws.Range("H78").Value = 3477.75
$ws.Range("J78").Value = 3477.75
$ws.Range("L78").Value = 31299.75
$ws.Range("N78").Value = -41283.75
$ws.Range("H81").Value = 1847.0625
$ws.Range("J81").Value = 2484.5454
$ws.Range("L81").Value = 7453.6362
$ws.Range("N81").Value = -9699.636200000001
$ws.Range("H84").Value = 1847.0625
$ws.Range("J84").Value = 2484.5454
$ws.Range("L84").Value = 22360.9086
$ws.Range("N84").Value = -33592.9086
$ws.Range("H131").Value = 1018.087
$ws.Range("I131").Value = 1285.1818
$ws.Range("J131").Value = 934.1429000000001
$ws.Range("K131").Value = 3855.5454
$ws.Range("L131").Value = 2802.4287
$ws.Range("M131").Value = 1184.4546
$ws.Range("N131").Value = -12882.4287
$ws.Range("H137").Value = 4596.3335
$ws.Range("I137").Value = 5266.6665
$ws.Range("J137").Value = 4462.2666
$ws.Range("K137").Value = 15799.9995
$ws.Range("L137").Value = 13386.7998
$ws.Range("M137").Value = -10699.9995
$ws.Range("N137").Value = -23586.7998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 39000
$ws.Range("J52").Value = 39000
$ws.Range("L52").Value = 39000
$ws.Range("N52").Value = -39518
$ws.Range("H70").Value = 5207.778
$ws.Range("I70").Value = 5601.25
$ws.Range("J70").Value = 5042.1055
$ws.Range("K70").Value = 5601.25
$ws.Range("L70").Value = 5042.1055
$ws.Range("M70").Value = -5331.25
$ws.Range("N70").Value = -5582.1055
$ws.Range("H73").Value = 5207.778
$ws.Range("I73").Value = 5601.25
$ws.Range("J73").Value = 5042.1055
$ws.Range("K73").Value = 5601.25
$ws.Range("L73").Value = 5042.1055
$ws.Range("M73").Value = -4665.25
$ws.Range("N73").Value = -6914.1055
$ws.Range("H80").Value = 31252788
$ws.Range("I80").Value = 35717044
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 35717044
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -35716046
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 31252788
$ws.Range("I83").Value = 35717044
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 178585220
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -178580228
$ws.Range("N83").Value = -24984
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -36134
$ws.Range("H113").Value = 1286.3
$ws.Range("I113").Value = 1305.5555
$ws.Range("J113").Value = 1113
$ws.Range("K113").Value = 1305.5555
$ws.Range("L113").Value = 1113
$ws.Range("M113").Value = 864.4445000000001
$ws.Range("N113").Value = -5453
$ws.Range("H132").Value = 3634.6667
$ws.Range("I132").Value = 2354.818
$ws.Range("J132").Value = 5042.5
$ws.Range("K132").Value = 7064.454000000001
$ws.Range("L132").Value = 15127.5
$ws.Range("M132").Value = -4534.454000000001
$ws.Range("N132").Value = -20187.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7445.2607
$ws.Range("I132").Value = 2741.4443
$ws.Range("K132").Value = 8224.332900000001
$ws.Range("M132").Value = -5694.332900000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 15865.728
$ws.Range("J54").Value = 15865.728
$ws.Range("L54").Value = 15865.728
$ws.Range("N54").Value = -16905.728
